# Fruta / hortaliza, semanal
# Insert a new data row above row 42 (pushing old rows 42-44 down to 43-45)
# and populate the new row 42 with the new weekly observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 42:44 down by inserting a new row at 42.
$ws.Rows.Item(42).Insert()

# Populate the freshly inserted row 42 with the new record.
$ws.Cells.Item(42, 1).Value = 10
$ws.Cells.Item(42, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(42, 3).Value = "La Araucanía"
$ws.Cells.Item(42, 4).Value = 44706
$ws.Cells.Item(42, 4).NumberFormat = $ws.Cells.Item(43, 4).NumberFormat
$ws.Cells.Item(42, 5).Value = 9
$ws.Cells.Item(42, 6).Value = "Fruta"
$ws.Cells.Item(42, 7).Value = 100107
$ws.Cells.Item(42, 8).Value = "Otros"
$ws.Cells.Item(42, 9).Value = 100107001
$ws.Cells.Item(42, 10).Value = "Caqui"
$ws.Cells.Item(42, 11).Value = "Mankaki"
$ws.Cells.Item(42, 12).Value = "Primera"
$ws.Cells.Item(42, 13).Value = 45
$ws.Cells.Item(42, 14).Value = 18000
$ws.Cells.Item(42, 15).Value = 18000
$ws.Cells.Item(42, 16).Value = 18000
$ws.Cells.Item(42, 17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(42, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(42, 19).Value = 1000
$ws.Cells.Item(42, 20).Value = 18
